$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet ships protected; temporarily unprotect so the cells below can be
# edited, then re-protect at the end so the sheet's protected state is
# restored (we can't reproduce the original legacy password hash, but we can
# restore the "sheet is protected" state).
$ws.Unprotect()

# Update the confidential disclaimer date (A12) from 2021-05-12 to 2021-05-13.
$ws.Range("A12").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-13 for illustrative purposes only and are subject to change."

# Update the recalculated Weight (D) and Percent Change (E) values for rows 2-9.
$ws.Range("D2").Value = 0.177392660699446
$ws.Range("E2").Value = 0.001808318264014286

$ws.Range("D3").Value = 0.1773906558095194
$ws.Range("E3").Value = 0.0009832841691248539

$ws.Range("D4").Value = 0.2249486497567567
$ws.Range("E4").Value = 0.001666666666666483

$ws.Range("D5").Value = 0.07998007139413028
$ws.Range("E5").Value = 0

$ws.Range("D6").Value = 0.07986679511328129

$ws.Range("D7").Value = 0.1202994102616281
$ws.Range("E7").Value = 0.0009832841691248539

$ws.Range("D8").Value = 0.1401217569652382
$ws.Range("E8").Value = 0

$ws.Range("E9").Value = 0.0009884107337796433

$ws.Protect("")
